# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" counts (column F) across the relevant sheets.

$wb = $excel.ActiveWorkbook

$wsExpo   = $wb.Worksheets.Item("展览")
$wsLocal  = $wb.Worksheets.Item("本地生活")
$wsAll    = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1) column F updates
$wsExpo.Range("F3").Value  = 7879
$wsExpo.Range("F8").Value  = 606
$wsExpo.Range("F12").Value = 864
$wsExpo.Range("F13").Value = 3148
$wsExpo.Range("F16").Value = 738
$wsExpo.Range("F19").Value = 457
$wsExpo.Range("F21").Value = 248
$wsExpo.Range("F23").Value = 314
$wsExpo.Range("F25").Value = 131
$wsExpo.Range("F26").Value = 107
$wsExpo.Range("F27").Value = 275
$wsExpo.Range("F32").Value = 520
$wsExpo.Range("F33").Value = 22
$wsExpo.Range("F34").Value = 33
$wsExpo.Range("F37").Value = 97

# 本地生活 (sheet3) column F update
$wsLocal.Range("F2").Value = 206

# 全部类型 (sheet4) column F updates
$wsAll.Range("F2").Value  = 206
$wsAll.Range("F5").Value  = 7879
$wsAll.Range("F10").Value = 606
$wsAll.Range("F14").Value = 864
$wsAll.Range("F16").Value = 3148
$wsAll.Range("F20").Value = 738
$wsAll.Range("F24").Value = 457
$wsAll.Range("F26").Value = 248
$wsAll.Range("F28").Value = 314
$wsAll.Range("F30").Value = 131
$wsAll.Range("F31").Value = 107
$wsAll.Range("F32").Value = 275
$wsAll.Range("F37").Value = 520
$wsAll.Range("F38").Value = 22
$wsAll.Range("F39").Value = 33
$wsAll.Range("F42").Value = 97
